$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 holds a list of "blog" card definitions referenced by shared
# strings. A new blog post (serial 160) was added to the front of the
# list, shifting the existing entries along:
#   B11: ser 159 -> ser 160 (new post)
#   D11: ser 156 -> ser 159 (previous B11 value)
#   I11: ser 154 -> ser 156 (previous D11 value)
# (the old I11 value, ser 154, drops off the list)

$oldD11 = $ws.Range("D11").Value()
$oldB11 = $ws.Range("B11").Value()

$ws.Range("I11").Value = $oldD11
$ws.Range("D11").Value = $oldB11
$ws.Range("B11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 160"
